$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: "CERTIFICADO N°{{n_certificado}}-2025" -> "...-2026"
# The certificate year suffix run changes its digit from 5 to 6.
# Scope the Find to the specific paragraph holding "CERTIFICADO N°"
# so the other "-2025" occurrence (N° DS {{num_ds}}-2025) is untouched.
# -----------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*CERTIFICADO N*") {
        $pr = $para.Range
        $found = $pr.Find.Execute("5", $true, $false, $false, $false, $false, $true, 1, $false, "6", 2)
        break
    }
}
Write-Output ("Certificado year updated: " + $found)

# -----------------------------------------------------------------
# Change 2: wrap the {{num_cara}} merge field with parentheses so the
# line reads "UNA ({{num_cara}}) CARA" instead of "UNA {{num_cara}} CARA".
# -----------------------------------------------------------------
$found2 = $d.Content.Find.Execute("{{num_cara}}", $true, $false, $false, $false, $false, $true, 1, $false, "({{num_cara}})", 2)
Write-Output ("num_cara wrapped in parentheses: " + $found2)
